$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Input Ports: "Port 2 : " -> "Port 4" + relocated _GoBack bookmark + " : "
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Port 2 : ", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "Port 4 : ", 2)

$pPort4 = $d.Paragraphs(4)
$splitPos = $pPort4.Range.Start + 6          # right after "Port 4"
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 2) Output Ports: renumber / reshuffle the three port paragraphs.
#    Old:  "Port 9 : " / "Port 14: " / "Port 15:"
#    New:  "Port 8: "  / "Port 9:"   / "Port 10 :"
#    Build the three replacement paragraphs from scratch (with real <w:tab/>
#    runs) right after the old block, then delete the old paragraphs.
# ---------------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newPort8  = '<w:p ' + $wNs + '>' + `
                 '<w:r><w:t>Port 8</w:t></w:r>' + `
                 '<w:r><w:t xml:space="preserve">: </w:t></w:r>' + `
                 '<w:r><w:tab/></w:r>' + `
                 '<w:r><w:tab/></w:r>' + `
                 '<w:r><w:t xml:space="preserve">Z1/Z8 </w:t></w:r>' + `
                 '<w:r><w:t>(Low Byte)</w:t></w:r>' + `
             '</w:p>'

$newPort9  = '<w:p ' + $wNs + '>' + `
                 '<w:r><w:t>Port 9</w:t></w:r>' + `
                 '<w:r><w:t>:</w:t></w:r>' + `
                 '<w:r><w:tab/></w:r>' + `
                 '<w:r><w:tab/><w:t xml:space="preserve">Z10/Z18 </w:t></w:r>' + `
                 '<w:r><w:t>(High Byte)</w:t></w:r>' + `
             '</w:p>'

$newPort10 = '<w:p ' + $wNs + '>' + `
                 '<w:r><w:t>Port 10 :</w:t></w:r>' + `
                 '<w:r><w:tab/></w:r>' + `
                 '<w:r><w:t>Keyboard reset strobe (P1.9)</w:t></w:r>' + `
             '</w:p>'

$pPort15 = $d.Paragraphs(8)
$insAt = $d.Range($pPort15.Range.End - 1, $pPort15.Range.End - 1)
$insAt.InsertXML($newPort8 + $newPort9 + $newPort10)

# Remove the three old paragraphs (still at indexes 6,7,8 - the new ones
# were inserted after them), highest index first so earlier indexes stay valid.
$d.Paragraphs(8).Range.Delete()
$d.Paragraphs(7).Range.Delete()
$d.Paragraphs(6).Range.Delete()

# ---------------------------------------------------------------------------
# 3) Remove the stray _GoBack bookmark left at the end of the document.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $endBm = $d.Bookmarks("_GoBack")
    if ($endBm.Range.Start -ne $splitPos) {
        $endBm.Delete()
    }
}

# The one we want is the one near "Port 4"; make sure only that one remains.
$count = 0
$last = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $last = $i
}
Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
